# daily auto push: 2026-01-18 22:34 UTC
# A new record for 2026/01/19 (Monday, hour slot "4", ranking 19) is
# inserted into the activity log right before the 2026/12/29 block,
# pushing that row and every row after it down by one (the table grows
# from 694 data rows to 695).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift row 653 and everything below it down by one row.
$ws.Rows.Item(653).EntireRow.Insert()

# Populate the newly inserted row with the new record.
# The leading apostrophe forces column A to stay plain text (matching
# the other date-like "yyyy/mm/dd" entries in the sheet) instead of
# Excel auto-converting it to a date serial number.
$ws.Cells.Item(653, 1).Value = "'2026/01/19"
$ws.Cells.Item(653, 2).Value = "月"
$ws.Cells.Item(653, 3).Value = 4
$ws.Cells.Item(653, 4).Value = 19
